# Powerpoint writer: consolidate text run nodes.
# Re-writes each title (and a couple of caption textboxes) so that a
# trailing single-space run gets folded into the preceding word run,
# instead of living on as its own <a:r> node. This is done by clearing
# the existing runs of the paragraph and re-inserting the consolidated
# chunks one at a time (each InsertAfter call produces its own <a:r>,
# while leaving the run immediately preceding it untouched).

function Set-ConsolidatedRuns {
    param(
        $TextRange,
        [string[]]$Chunks
    )
    $TextRange.Delete()
    foreach ($chunk in $Chunks) {
        $TextRange.InsertAfter($chunk) | Out-Null
    }
}

$p = $ppt.ActivePresentation

# Slide 1: "Slide 1 (Content)"
$s = $p.Slides.Item(1)
Set-ConsolidatedRuns $s.Shapes.Item(1).TextFrame.TextRange @("Slide ", "1 ", "(Content)")

# Slide 2: "Slide 2 (Content)"
$s = $p.Slides.Item(2)
Set-ConsolidatedRuns $s.Shapes.Item(1).TextFrame.TextRange @("Slide ", "2 ", "(Content)")

# Slide 3: "Slide 3 (Content)"
$s = $p.Slides.Item(3)
Set-ConsolidatedRuns $s.Shapes.Item(1).TextFrame.TextRange @("Slide ", "3 ", "(Content)")

# Slide 4: "Slide 4 (Content)"
$s = $p.Slides.Item(4)
Set-ConsolidatedRuns $s.Shapes.Item(1).TextFrame.TextRange @("Slide ", "4 ", "(Content)")

# Slide 5: "Slide 5 (Two Content)"
$s = $p.Slides.Item(5)
Set-ConsolidatedRuns $s.Shapes.Item(1).TextFrame.TextRange @("Slide ", "5 ", "(Two ", "Content)")

# Slide 6: "Slide 6 (Two Content Right)" + "an image" caption
$s = $p.Slides.Item(6)
Set-ConsolidatedRuns $s.Shapes.Item(1).TextFrame.TextRange @("Slide ", "6 ", "(Two ", "Content ", "Right)")
Set-ConsolidatedRuns $s.Shapes.Item(3).TextFrame.TextRange @("an ", "image")

# Slide 7: "Slide 7 (Content with Caption)" + "An image" caption
$s = $p.Slides.Item(7)
Set-ConsolidatedRuns $s.Shapes.Item(1).TextFrame.TextRange @("Slide ", "7 ", "(Content ", "with ", "Caption)")
Set-ConsolidatedRuns $s.Shapes.Item(4).TextFrame.TextRange @("An ", "image")

# Slide 8: "Slide 8 (Comparison)" + "An image" caption
$s = $p.Slides.Item(8)
Set-ConsolidatedRuns $s.Shapes.Item(1).TextFrame.TextRange @("Slide ", "8 ", "(Comparison)")
Set-ConsolidatedRuns $s.Shapes.Item(4).TextFrame.TextRange @("An ", "image")

# Slide 9: "Slide 10 (Content)"
$s = $p.Slides.Item(9)
Set-ConsolidatedRuns $s.Shapes.Item(1).TextFrame.TextRange @("Slide ", "10 ", "(Content)")

# Slide 10: "Slide 11 (Content)"
$s = $p.Slides.Item(10)
Set-ConsolidatedRuns $s.Shapes.Item(1).TextFrame.TextRange @("Slide ", "11 ", "(Content)")

# Slide 11: "Slide 12 (Content)"
$s = $p.Slides.Item(11)
Set-ConsolidatedRuns $s.Shapes.Item(1).TextFrame.TextRange @("Slide ", "12 ", "(Content)")
